$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.180128812789917
$ws.Range("B1").Value = 2.388216733932495
$ws.Range("C1").Value = 3.589848041534424
$ws.Range("D1").Value = 1.971291780471802
$ws.Range("E1").Value = 1.206217288970947
